$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price values so they are not
# auto-converted to numbers by Excel (source data are text strings).
$textCells = @("D4", "D5", "D8", "D11", "D13", "D15", "D17", "D19", "D22", "D25", "D26", "D29", "D30", "D33", "D40", "D42", "D44", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = '34.110.72'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '1.778.65'
$ws.Range("E3").Value = '  -2.58%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = '225.24'
$ws.Range("E5").Value = '  -1.94%  '
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").Value = '31.59'
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("E9").Value = '  -1.05%  '
$ws.Range("E10").Value = '  -2.36%  '
$ws.Range("D11").Value = '0.0930'
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").Value = '2.034.75'
$ws.Range("E12").Value = '  -2.57%  '
$ws.Range("D13").Value = '11.06'
$ws.Range("E13").Value = '  +6.41%  '
$ws.Range("D14").Value = '1.767.10'
$ws.Range("E14").Value = '  -3.30%  '
$ws.Range("D15").Value = '0.624'
$ws.Range("E15").Value = '  -3.24%  '
$ws.Range("D16").Value = '34.108.63'
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("D17").Value = '4.20'
$ws.Range("E17").Value = '  -1.89%  '
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("D19").Value = '254.15'
$ws.Range("E19").Value = '  -2.01%  '
$ws.Range("D20").Value = '0.0₃0736'
$ws.Range("E20").Value = '  -2.27%  '
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("D22").Value = '10.34'
$ws.Range("E22").Value = '  -2.59%  '
$ws.Range("E24").Value = '  -3.32%  '
$ws.Range("D25").Value = '156.74'
$ws.Range("E25").Value = '  -1.30%  '
$ws.Range("D26").Value = '16.38'
$ws.Range("E26").Value = '  -1.91%  '
$ws.Range("E27").Value = '  -2.36%  '
$ws.Range("E28").Value = '  -1.47%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("D30").Value = '3.76'
$ws.Range("E31").Value = '  -0.78%  '
$ws.Range("E32").Value = '  -1.71%  '
$ws.Range("D33").Value = '3.58'
$ws.Range("E33").Value = '  +0.60%  '
$ws.Range("E34").Value = '  +1.25%  '
$ws.Range("D35").Value = '1.440.21'
$ws.Range("E35").Value = '  -7.25%  '
$ws.Range("E36").Value = '  -3.92%  '
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("E38").Value = '  -1.54%  '
$ws.Range("E39").Value = '  +0.67%  '
$ws.Range("D40").Value = '82.81'
$ws.Range("E41").Value = '  +0.34%  '
$ws.Range("D42").Value = '0.886'
$ws.Range("E42").Value = '  -3.86%  '
$ws.Range("E43").Value = '  -5.79%  '
$ws.Range("D44").Value = '0.0512'
$ws.Range("E44").Value = '  -2.80%  '
$ws.Range("E45").Value = '  -2.05%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '5.81'
$ws.Range("E46").Value = '  +1.07%  '
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '1.934.31'
$ws.Range("E47").Value = '  -2.71%  '
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").Value = '0.999'
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '11.97'
$ws.Range("E49").Value = '  -1.00%  '
$ws.Range("D50").Value = '98.35'
$ws.Range("E50").Value = '  +1.00%  '
$ws.Range("D51").Value = '49.47'
$ws.Range("E51").Value = '  -7.02%  '
